$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 <- old row 11 data (111702486)
$ws.Range("A8").Value = 111702486
$ws.Range("B8").Value = 90678
$ws.Range("E8").Value = 4366
$ws.Range("F8").Value = "Skarp dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum peckii"
$ws.Range("H8").Value = "Banker"
$ws.Range("P8").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q8").Value = 517080.8398438052
$ws.Range("R8").Value = 6574959.907818918

# Row 9 <- old row 8 data (111702506)
$ws.Range("A9").Value = 111702506
$ws.Range("B9").Value = 90687
$ws.Range("E9").Value = 5964
$ws.Range("F9").Value = "Fjällig taggsvamp s.str."
$ws.Range("G9").Value = "Sarcodon imbricatus s.str."
$ws.Range("H9").Value = "(L.:Fr.) P.Karst."
$ws.Range("P9").Value = "Kyrkogården, Nrk"
$ws.Range("Q9").Value = 517093.6249861007
$ws.Range("R9").Value = 6574959.965416327

# Row 10 <- old row 12 data (111702420)
$ws.Range("A10").Value = 111702420
$ws.Range("B10").Value = 90709
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 5448
$ws.Range("F10").Value = "Svartvit taggsvamp"
$ws.Range("G10").Value = "Phellodon connatus"
$ws.Range("H10").Value = "(Schultz) nom.prov"
$ws.Range("I10").Value = "1"
$ws.Range("J10").Value = "fruktkroppar"
$ws.Range("P10").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q10").Value = 517086.1792710476
$ws.Range("R10").Value = 6574909.900584662

# Row 11 <- old row 10 data (111702393)
$ws.Range("A11").Value = 111702393
$ws.Range("B11").Value = 89183
$ws.Range("E11").Value = 3215
$ws.Range("F11").Value = "Rödgul trumpetsvamp"
$ws.Range("G11").Value = "Craterellus lutescens"
$ws.Range("H11").Value = "(Fr.) Fr."
$ws.Range("P11").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q11").Value = 517070.2129045375
$ws.Range("R11").Value = 6574934.844418272

# Row 12 <- old row 9 data (111702400)
$ws.Range("A12").Value = 111702400
$ws.Range("B12").Value = 90687
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 5964
$ws.Range("F12").Value = "Fjällig taggsvamp s.str."
$ws.Range("G12").Value = "Sarcodon imbricatus s.str."
$ws.Range("H12").Value = "(L.:Fr.) P.Karst."
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("P12").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q12").Value = 517073.2951468225
$ws.Range("R12").Value = 6574931.795150192
